# Update the product URLs on Sheet1 to point to FairPrice search results
# instead of the (now dead/changed) direct product-page links.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value  = "https://www.fairprice.com.sg/search?query=Fiber%20Jelly%20Grape"
$ws.Range("D4").Value  = "https://www.fairprice.com.sg/search?query=Aw%27s%20Market%20Fresh%20Indonesian%20Pork%20Belly%20(Sliced)"
$ws.Range("D6").Value  = "https://www.fairprice.com.sg/search?query=Pokka%20Ice%20Lemon%20Tea"
$ws.Range("D8").Value  = "https://www.fairprice.com.sg/search?query=Marigold%20HL%20Milk%20-%20Chocolate"
$ws.Range("D10").Value = "https://www.fairprice.com.sg/search?query=Milo%20Instant%20Chocalate%20Malt%20Drink"
$ws.Range("D11").Value = "https://www.fairprice.com.sg/search?query=Maggi%20Big%20Curry%20Noodle"
$ws.Range("D12").Value = "https://www.fairprice.com.sg/search?query=Lipton%20Yellow%20Label%20Tea"
$ws.Range("D13").Value = "https://www.fairprice.com.sg/search?query=Panteen%20Shampoo%20Hail%20Fall%20Control"
